$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom")

# Row 9: Calc moved from B to C, replaced by F13 (Obsidian)
$ws.Range("B9").Value = "F13 (Obsidian)"
$ws.Range("C9").Value = "Calc"

# Row 3: AltGr (RAlt) -> MO(3)
$ws.Range("B3").Value = "MO(3)"

# New rows 14-18: Umlauts / Vim arrows distributed onto Layer 3 (Fn)
$ws.Range("A14").Value = "ß / -"
$ws.Range("C14").Value = "F17"
$ws.Range("C15").Value = "F18"
$ws.Range("C16").Value = "F19"
$ws.Range("C17").Value = "F20"
$ws.Range("C18").Value = "F21"

$ws.Range("A15").Value = "ü / ["
$ws.Range("A16").Value = "ö / ;"
$ws.Range("A17").Value = "ä / '"
$ws.Range("A18").Value = "€ / 5"

# Update selection to match new active cell after edits
$ws.Range("A19").Select()

$ws.Range("A14:C18").SetPhonetic()
